# Apply updated cryptocurrency price/volume data.
# Each target cell is written via a literal-text formula (="...")
# then flattened to a static value with Copy + PasteSpecial(xlPasteValues,-4163).
# This guarantees the text is stored exactly as given (no numeric/date
# auto-conversion of values like "1.00" or "0.100") while leaving cell
# styles/number-formats completely untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=`"27.555.34`""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = "=`"  -0.56%  `""
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("D3").Formula = "=`"1.623.68`""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = "=`"  -1.33%  `""
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)
$ws.Range("E4").Formula = "=`"  -0.05%  `""
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("D5").Formula = "=`"211.45`""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = "=`"  -0.81%  `""
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=`"0.526`""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Formula = "=`"  -0.94%  `""
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)
$ws.Range("E7").Formula = "=`"  -0.05%  `""
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)
$ws.Range("D8").Formula = "=`"23.16`""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Formula = "=`"  -0.48%  `""
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("E9").Formula = "=`"  +1.91%  `""
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("E10").Formula = "=`"  -0.28%  `""
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("E11").Formula = "=`"  -0.42%  `""
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)
$ws.Range("D12").Formula = "=`"1.854.05`""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = "=`"  -1.32%  `""
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$ws.Range("D13").Formula = "=`"1.623.27`""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = "=`"  -1.29%  `""
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("E14").Formula = "=`"  +0.23%  `""
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)
$ws.Range("E15").Formula = "=`"  -1.94%  `""
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("D16").Formula = "=`"65.26`""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = "=`"  +0.71%  `""
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("D17").Formula = "=`"27.517.49`""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Formula = "=`"  -0.65%  `""
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)
$ws.Range("D18").Formula = "=`"231.63`""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Formula = "=`"  -0.11%  `""
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("E19").Formula = "=`"  -0.77%  `""
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)
$ws.Range("D20").Formula = "=`"7.55`""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = "=`"  -1.06%  `""
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("E21").Formula = "=`"  -0.07%  `""
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)
$ws.Range("D22").Formula = "=`"10.44`""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Formula = "=`"  +2.85%  `""
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("E23").Formula = "=`"  +0.83%  `""
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("D24").Formula = "=`"2.08`""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Formula = "=`"  +6.39%  `""
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)
$ws.Range("D25").Formula = "=`"150.37`""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = "=`"  +0.17%  `""
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)
$ws.Range("E26").Formula = "=`"  -0.71%  `""
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)
$ws.Range("E27").Formula = "=`"  -0.65%  `""
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("B28").Formula = "=`"BinanceUSD`""
$ws.Range("B28").Copy()
$ws.Range("B28").PasteSpecial(-4163)
$ws.Range("C28").Formula = "=`"https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd`""
$ws.Range("C28").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=`"1.00`""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Formula = "=`"  +0.00%  `""
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("B29").Formula = "=`"EthereumClassic`""
$ws.Range("B29").Copy()
$ws.Range("B29").PasteSpecial(-4163)
$ws.Range("C29").Formula = "=`"https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc`""
$ws.Range("C29").Copy()
$ws.Range("C29").PasteSpecial(-4163)
$ws.Range("D29").Formula = "=`"15.55`""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Formula = "=`"  -0.71%  `""
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("E30").Formula = "=`"  -0.97%  `""
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("E31").Formula = "=`"  -0.60%  `""
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("D32").Formula = "=`"3.27`""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Formula = "=`"  -0.74%  `""
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=`"1.466.37`""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Formula = "=`"  +1.82%  `""
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)
$ws.Range("D34").Formula = "=`"3.07`""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Formula = "=`"  -2.23%  `""
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)
$ws.Range("E35").Formula = "=`"  -2.85%  `""
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)
$ws.Range("D36").Formula = "=`"2.35`""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Formula = "=`"  +0.07%  `""
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)
$ws.Range("D37").Formula = "=`"0.949`""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Formula = "=`"  +6.99%  `""
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)
$ws.Range("E38").Formula = "=`"  +0.64%  `""
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)
$ws.Range("B39").Formula = "=`"ARBITRUM`""
$ws.Range("B39").Copy()
$ws.Range("B39").PasteSpecial(-4163)
$ws.Range("C39").Formula = "=`"https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb`""
$ws.Range("C39").Copy()
$ws.Range("C39").PasteSpecial(-4163)
$ws.Range("D39").Formula = "=`"0.871`""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Formula = "=`"  -0.86%  `""
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)
$ws.Range("B40").Formula = "=`"ImmutableX`""
$ws.Range("B40").Copy()
$ws.Range("B40").PasteSpecial(-4163)
$ws.Range("C40").Formula = "=`"https://coinranking.com/coin/Z96jIvLU7+immutablex-imx`""
$ws.Range("C40").Copy()
$ws.Range("C40").PasteSpecial(-4163)
$ws.Range("D40").Formula = "=`"0.555`""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = "=`"  -2.74%  `""
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)
$ws.Range("E41").Formula = "=`"  -0.03%  `""
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)
$ws.Range("E42").Formula = "=`"  -2.12%  `""
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)
$ws.Range("D43").Formula = "=`"67.63`""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = "=`"  +0.28%  `""
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)
$ws.Range("B44").Formula = "=`"mCoin`""
$ws.Range("B44").Copy()
$ws.Range("B44").PasteSpecial(-4163)
$ws.Range("C44").Formula = "=`"https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin`""
$ws.Range("C44").Copy()
$ws.Range("C44").PasteSpecial(-4163)
$ws.Range("D44").Formula = "=`"2.45`""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = "=`"  -1.14%  `""
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)
$ws.Range("B45").Formula = "=`"MXToken`""
$ws.Range("B45").Copy()
$ws.Range("B45").PasteSpecial(-4163)
$ws.Range("C45").Formula = "=`"https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx`""
$ws.Range("C45").Copy()
$ws.Range("C45").PasteSpecial(-4163)
$ws.Range("D45").Formula = "=`"2.20`""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Formula = "=`"  -2.16%  `""
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)
$ws.Range("B46").Formula = "=`"FraxShare`""
$ws.Range("B46").Copy()
$ws.Range("B46").PasteSpecial(-4163)
$ws.Range("C46").Formula = "=`"https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs`""
$ws.Range("C46").Copy()
$ws.Range("C46").PasteSpecial(-4163)
$ws.Range("D46").Formula = "=`"5.30`""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Formula = "=`"  -5.13%  `""
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)
$ws.Range("B47").Formula = "=`"RenderToken`""
$ws.Range("B47").Copy()
$ws.Range("B47").PasteSpecial(-4163)
$ws.Range("C47").Formula = "=`"https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr`""
$ws.Range("C47").Copy()
$ws.Range("C47").PasteSpecial(-4163)
$ws.Range("D47").Formula = "=`"1.75`""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = "=`"  +0.55%  `""
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)
$ws.Range("B48").Formula = "=`"RocketPoolETH`""
$ws.Range("B48").Copy()
$ws.Range("B48").PasteSpecial(-4163)
$ws.Range("C48").Formula = "=`"https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth`""
$ws.Range("C48").Copy()
$ws.Range("C48").PasteSpecial(-4163)
$ws.Range("D48").Formula = "=`"1.764.60`""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = "=`"  -1.33%  `""
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)
$ws.Range("B49").Formula = "=`"Quant`""
$ws.Range("B49").Copy()
$ws.Range("B49").PasteSpecial(-4163)
$ws.Range("C49").Formula = "=`"https://coinranking.com/coin/bauj_21eYVwso+quant-qnt`""
$ws.Range("C49").Copy()
$ws.Range("C49").PasteSpecial(-4163)
$ws.Range("D49").Formula = "=`"87.32`""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Formula = "=`"  +1.98%  `""
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)
$ws.Range("B50").Formula = "=`"BabyDogeCoin`""
$ws.Range("B50").Copy()
$ws.Range("B50").PasteSpecial(-4163)
$ws.Range("C50").Formula = "=`"https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge`""
$ws.Range("C50").Copy()
$ws.Range("C50").PasteSpecial(-4163)
$ws.Range("D50").Formula = "=`"0.0₆0105`""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Formula = "=`"  -2.01%  `""
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)
$ws.Range("B51").Formula = "=`"Algorand`""
$ws.Range("B51").Copy()
$ws.Range("B51").PasteSpecial(-4163)
$ws.Range("C51").Formula = "=`"https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo`""
$ws.Range("C51").Copy()
$ws.Range("C51").PasteSpecial(-4163)
$ws.Range("D51").Formula = "=`"0.100`""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Formula = "=`"  +1.49%  `""
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)

$excel.CutCopyMode = 0
